# v 4.2.2 modifications to renal + add rep hfh
#
# Column A of "Sheet 1" holds an ordered list of variable names (one per
# row). This script applies the following logical changes:
#   1. Reorder "sos_com_dialysis" so it comes right after "sos_com_renal"
#      (before sos_com_hyperkalemia / sos_com_hypokalemia).
#   2. Rename "sos_out_hosprenal" -> "sos_out_hosprenalacute" and
#      "sos_outtime_hosprenal" -> "sos_outtime_hosprenalacute".
#   3. Insert two new variables right after those renamed rows:
#      "sos_out_renalendstage" and "sos_outtime_renalendstage".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Reorder sos_com_dialysis / sos_com_hyperkalemia / sos_com_hypokalemia ---
# Before (rows 165-167): sos_com_hyperkalemia, sos_com_hypokalemia, sos_com_dialysis
# After  (rows 165-167): sos_com_dialysis, sos_com_hyperkalemia, sos_com_hypokalemia
$ws.Range("A165").Value = "sos_com_dialysis"
$ws.Range("A166").Value = "sos_com_hyperkalemia"
$ws.Range("A167").Value = "sos_com_hypokalemia"

# --- 2. Rename sos_out_hosprenal / sos_outtime_hosprenal ---
$ws.Range("A203").Value = "sos_out_hosprenalacute"
$ws.Range("A204").Value = "sos_outtime_hosprenalacute"

# --- 3. Insert two new rows for the new "renal end stage" outcome variables ---
$ws.Rows("205:206").Insert()
$ws.Range("A205").Value = "sos_out_renalendstage"
$ws.Range("A206").Value = "sos_outtime_renalendstage"
